$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.142.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.246.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.64%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "90.62"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.21%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.48"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +9.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.63"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0791"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.10%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.55"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.593.48"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.05"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.195.44"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.748"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.096.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.80"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0899"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.84"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.62"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "239.61"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.41%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.85"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.68"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.37%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.85"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.05"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.32%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.15"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0730"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.21%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.115"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.32"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.074.66"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.01"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +14.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0276"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +11.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.468.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.50"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.28%  "
